$wb = $excel.ActiveWorkbook

# Sheets "展览" and "全部类型" contain identical data tables that both need updating.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("G2").Value = 55
    $ws.Range("F4").Value = 31
}
